# Update forest data - 2025-11-24 12:19
#
# Sheet "Previously added": append the row that used to be the sole data
# row on sheet "New" (it becomes row 272).
# Sheet "New": replace its row 2 with a fresh listing and append 7 more
# fresh listings (rows 3-9).

$wb = $excel.ActiveWorkbook
$wsPrev = $wb.Worksheets.Item("Previously added")
$wsNew  = $wb.Worksheets.Item("New")

# Untouched, known-good formatting template (link / text / text / text /
# text / date) used to restore each new/changed row's look after the
# value + hyperlink writes (which can reset a cell's style).
$wsPrev.Range("A271:F271").Copy()

# ---------------------------------------------------------------------
# 1) "Previously added": add row 272 with the data that (pre-edit)
#    lived on "New"!A2:F2.
# ---------------------------------------------------------------------
$wsPrev.Range("A272").Value = "https://www.ss.com/msg/lv/real-estate/wood/talsi-and-reg/dundagas-pag/mjdlj.html"
$wsPrev.Range("B272").Value = "169 000 €"
$wsPrev.Range("C272").Value = "Talsi un raj."
$wsPrev.Range("D272").Value = "59 ha."
$wsPrev.Range("E272").Value = "'88500020198"
$wsPrev.Range("F272").Value = 45982.368055555555
$wsPrev.Hyperlinks.Add($wsPrev.Range("A272"), "https://www.ss.com/msg/lv/real-estate/wood/talsi-and-reg/dundagas-pag/mjdlj.html")

$wsPrev.Range("A272:F272").PasteSpecial(-4122)  # xlPasteFormats, restores s="3,4,4,4,4,2"

# ---------------------------------------------------------------------
# 2) "New": overwrite row 2 with a new listing, then append rows 3-9.
# ---------------------------------------------------------------------
$newRows = @(
    @{ Row = 2; A = "https://www.ss.com/msg/lv/real-estate/wood/aluksne-and-reg/jaunlaicenes-pag/kxjph.html"; B = "692 000 €";     C = "Alūksne un raj.";   D = "30 ha.";   E = "36600010045"; F = 45984.62013888889 },
    @{ Row = 3; A = "https://www.ss.com/msg/lv/real-estate/wood/balvi-and-reg/berzpils-pag/hiidj.html?_gl=1*1u7cuf1*_up*MQ..*_ga*ODQ5MTg2ODA0LjE3NjM5ODY3Mzk.*_ga_ZCGHC71BQ2*czE3NjM5ODY3MzkkbzEkZzAkdDE3NjM5ODY3MzkkajYwJGwwJGgw"; B = "24 500 €";     C = "Balvi un raj.";    D = "1.80 ha."; E = "38500050070"; F = 45982.84305555555 },
    @{ Row = 4; A = "https://www.ss.com/msg/lv/real-estate/wood/balvi-and-reg/berzpils-pag/hbxpx.html?_gl=1*17t3g13*_up*MQ..*_ga*MTI2Mzg3NDExOS4xNzYzOTg2NzQw*_ga_ZCGHC71BQ2*czE3NjM5ODY3MzkkbzEkZzAkdDE3NjM5ODY3MzkkajYwJGwwJGgw"; B = "65 000 €"; C = "Balvi un raj.";    D = "6 ha.";    E = "38500030175"; F = 45982.70277777778 },
    @{ Row = 5; A = "https://www.ss.com/msg/lv/real-estate/wood/gulbene-and-reg/litenes-pag/gmxfh.html";       B = "130 000 €"; C = "Gulbene un raj.";  D = "20 ha.";   E = "50680070035"; F = 45984.84722222222 },
    @{ Row = 6; A = "https://www.ss.com/msg/lv/real-estate/wood/jekabpils-and-reg/atasienes-pag/cndok.html";  B = "30 000 €";  C = "Jēkabpils un raj."; D = "10 ha.";   E = "56460030122"; F = 45983.822222222225 },
    @{ Row = 7; A = "https://www.ss.com/msg/lv/real-estate/wood/kraslava-and-reg/aulejas-pag/opgnm.html";     B = "8 000 €";   C = "Krāslava un raj."; D = "1 ha.";    E = "60480030114"; F = 45985.59583333333 },
    @{ Row = 8; A = "https://www.ss.com/msg/lv/real-estate/wood/ludza-and-reg/malnavas-pag/iiemg.html";       B = "12 000 €";  C = "Ludza un raj.";   D = "1 ha.";    E = "68680060098"; F = 45983.424305555556 },
    @{ Row = 9; A = "https://www.ss.com/msg/lv/real-estate/wood/talsi-and-reg/ives-pag/bxkod.html";           B = "1 080 000 €"; C = "Talsi un raj."; D = "40 ha.";   E = "88580020037"; F = 45983.527083333334 }
)

# Row 2 carries a hyperlink pointing at the old listing; drop it before
# the cell gets new content.
$wsNew.Range("A2").Hyperlinks.Delete()

foreach ($r in $newRows) {
    $wsNew.Range("A$($r.Row)").Value = $r.A
    $wsNew.Range("B$($r.Row)").Value = $r.B
    $wsNew.Range("C$($r.Row)").Value = $r.C
    $wsNew.Range("D$($r.Row)").Value = $r.D
    $wsNew.Range("E$($r.Row)").Value = "'" + $r.E
    $wsNew.Range("F$($r.Row)").Value = $r.F

    $wsNew.Hyperlinks.Add($wsNew.Range("A$($r.Row)"), $r.A)

    $wsNew.Range("A$($r.Row):F$($r.Row)").PasteSpecial(-4122)  # xlPasteFormats
}

Write-Host "Done."
